# Update "Schedule" sheet (sheet1): a new optimisation run inserts an extra
# scheduled block, shifting the old row 4 down to row 5, and refreshes the
# numeric results for every block.
$wb = $excel.ActiveWorkbook
$schedule = $wb.Worksheets.Item("Schedule")
$detail = $wb.Worksheets.Item("Detailed")

# Insert a new row 4 (pushes the old row 4 down to row 5, carrying formatting).
$schedule.Rows.Item(4).Insert()

$scheduleRows = [ordered]@{
    2 = @{ B=46056.27083333334; C=6.5; D=24.57; E=865.85679375; F=35.24040674603174 }
    3 = @{ A=46056.29166666666; B=46056.66666666666; C=9; D=34.02; E=728.5141402499999; F=21.41428983686067 }
    4 = @{ A=46056.97916666666; B=46057.14583333334; C=4; D=15.12; E=759.9122115; F=50.2587441468254 }
    5 = @{ A=46057.3125; B=46057.66666666666; C=8.5; D=32.13; E=143.36233275; F=4.461946241830066 }
}

foreach ($r in $scheduleRows.Keys) {
    $row = $scheduleRows[$r]
    foreach ($col in $row.Keys) {
        $schedule.Range("$col$r").Value = $row[$col]
    }
}

# Update "Detailed" sheet (sheet2): refreshed Price forecast/history values,
# a few Type (forecast -> historical) reclassifications, and a handful of
# Pump_Status flips, per the latest optimisation run.
$detailRows = [ordered]@{
    15 = @{ E="OFF" }
    38 = @{ B=101.25 }
    39 = @{ B=84.79000000000001 }
    40 = @{ B=104.91701; C="historical" }
    41 = @{ B=105.79; C="historical" }
    42 = @{ B=103.03022; C="historical" }
    43 = @{ B=106.12755; C="historical" }
    44 = @{ B=108.89; C="historical" }
    45 = @{ B=108.89; C="historical" }
    46 = @{ B=101.33; C="historical" }
    47 = @{ B=101.53606; C="historical" }
    48 = @{ B=100.91363; C="historical" }
    49 = @{ B=84.79000000000001; E="ON" }
    50 = @{ B=97.73156; E="ON" }
    51 = @{ B=101.25; E="ON" }
    52 = @{ B=102.50664 }
    53 = @{ B=101.25 }
    54 = @{ B=101.25 }
    55 = @{ B=101.25 }
    56 = @{ B=89.36893999999999 }
    57 = @{ B=84.90018999999999; E="OFF" }
    58 = @{ E="OFF" }
    59 = @{ B=100.3; E="OFF" }
    60 = @{ B=105.79 }
    61 = @{ B=110.39857 }
    62 = @{ B=108.89 }
    63 = @{ B=115.90225 }
    64 = @{ B=84.79000000000001 }
    65 = @{ B=57.06007; E="ON" }
    66 = @{ B=38.68104 }
    67 = @{ B=32.96129 }
    68 = @{ B=21.5418 }
    69 = @{ B=0.7 }
    70 = @{ B=0.01101 }
    71 = @{ B=0 }
    72 = @{ B=-2.54301 }
    73 = @{ B=-1.15096 }
    74 = @{ B=0 }
    75 = @{ B=-4.67865 }
    76 = @{ B=-4.70553 }
    78 = @{ B=-1.04393 }
    79 = @{ B=0.00002 }
    80 = @{ B=0.51 }
    81 = @{ B=9.185280000000001 }
    82 = @{ B=20.89615 }
    83 = @{ B=36.0601 }
    84 = @{ B=54.25499 }
    85 = @{ B=63.05162 }
    86 = @{ B=69.26626 }
    87 = @{ B=102.96323 }
    88 = @{ B=135.5606 }
    89 = @{ B=141.98066 }
    90 = @{ B=169.53226 }
    91 = @{ B=150.50671 }
    92 = @{ B=142.08175 }
    93 = @{ B=118.30564 }
    94 = @{ B=115.04338 }
    95 = @{ B=113.3996 }
    96 = @{ B=111.84804 }
    97 = @{ B=115.71339 }
}

foreach ($r in $detailRows.Keys) {
    $row = $detailRows[$r]
    foreach ($col in $row.Keys) {
        $detail.Range("$col$r").Value = $row[$col]
    }
}
